# Insert a new data row at row 51 (pushes existing rows 51-120 down to 52-121)
# and populate it with the new weekly price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("51:51").Insert()

$ws.Range("A51").Value = 4
$ws.Range("B51").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C51").Value = "Los Lagos"
$ws.Range("D51").Value = 44971
$ws.Range("E51").Value = 10
$ws.Range("F51").Value = 100112043
$ws.Range("G51").Value = "Pepino dulce"
$ws.Range("H51").Value = "Cultivar IV Región"
$ws.Range("I51").Value = "Especial"
$ws.Range("J51").Value = 50
$ws.Range("K51").Value = 24000
$ws.Range("L51").Value = 24000
$ws.Range("M51").Value = 24000
$ws.Range("N51").Value = "`$/bandeja 18 kilos"
$ws.Range("O51").Value = "Provincia de Limarí"
$ws.Range("P51").Value = 1333
$ws.Range("Q51").Value = 18
$ws.Range("R51").Value = "Hortaliza"
